$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "30.068.29"
$ws.Range("E2").Value = "  +0.10%  "
# Row 3
$ws.Range("D3").Value = "1.882.12"
$ws.Range("E3").Value = "  +0.40%  "
# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9986"
$ws.Range("E4").Value = "  -0.16%  "
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "243.67"
$ws.Range("E5").Value = "  -1.87%  "
# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9985"
$ws.Range("E6").Value = "  -0.17%  "
# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4962"
$ws.Range("E7").Value = "  -0.03%  "
# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "44.45"
$ws.Range("E8").Value = "  -2.77%  "
# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2909"
# Row 10
$ws.Range("E10").Value = "  +0.89%  "
# Row 11
$ws.Range("D11").Value = "1.879.14"
$ws.Range("E11").Value = "  +0.35%  "
# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "16.83"
$ws.Range("E12").Value = "  -1.37%  "
# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.07189"
$ws.Range("E13").Value = "  +0.08%  "
# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6640"
$ws.Range("E14").Value = "  +0.50%  "
# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "85.67"
$ws.Range("E15").Value = "  +0.67%  "
# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "4.838"
$ws.Range("E16").Value = "  +0.96%  "
# Row 17
$ws.Range("D17").Value = "30.051.43"
$ws.Range("E17").Value = "  +0.15%  "
# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000007795"
$ws.Range("E18").Value = "  +4.06%  "
# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.9982"
$ws.Range("E19").Value = "  -0.30%  "
# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.77"
$ws.Range("E20").Value = "  -0.42%  "
# Row 21
$ws.Range("D21").Value = "2.122.28"
$ws.Range("E21").Value = "  +0.43%  "
# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9981"
$ws.Range("E22").Value = "  -0.19%  "
# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.759"
$ws.Range("E23").Value = "  +0.66%  "
# Row 24
$ws.Range("B24").Value = "Cosmos"
$ws.Range("C24").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.157"
$ws.Range("E24").Value = "  +1.50%  "
# Row 25
$ws.Range("B25").Value = "Chainlink"
$ws.Range("C25").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "5.600"
$ws.Range("E25").Value = "  +1.86%  "
# Row 26
$ws.Range("E26").Value = "  +4.46%  "
# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "136.22"
$ws.Range("E27").Value = "  +0.76%  "
# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "16.79"
$ws.Range("E28").Value = "  +0.57%  "
# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.909"
$ws.Range("E29").Value = "  -2.29%  "
# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.380"
$ws.Range("E30").Value = "  -0.43%  "
# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.161"
$ws.Range("E31").Value = "  -1.01%  "
# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.08686"
$ws.Range("E32").Value = "  +1.23%  "
# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.944"
$ws.Range("E33").Value = "  +1.69%  "
# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.04995"
$ws.Range("E34").Value = "  -1.31%  "
# Row 35
$ws.Range("E35").Value = "  -2.65%  "
# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7031"
$ws.Range("E36").Value = "  +3.10%  "
# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.657"
$ws.Range("E37").Value = "  -1.68%  "
# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.702"
$ws.Range("E38").Value = "  -1.42%  "
# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.198"
$ws.Range("E39").Value = "  -4.80%  "
# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.9334"
$ws.Range("E40").Value = "  -3.01%  "
# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.01643"
$ws.Range("E41").Value = "  +1.10%  "
# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.952"
$ws.Range("E42").Value = "  -1.73%  "
# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.9990"
$ws.Range("E43").Value = "  -0.14%  "
# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.4189"
$ws.Range("E44").Value = "  +0.04%  "
# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "101.38"
$ws.Range("E45").Value = "  -1.69%  "
# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "7.511"
$ws.Range("E46").Value = "  +0.86%  "
# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.1260"
$ws.Range("E47").Value = "  +0.49%  "
# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.05716"
$ws.Range("E48").Value = "  +1.67%  "
# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "32.35"
$ws.Range("E49").Value = "  -0.04%  "
# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.252"
$ws.Range("E50").Value = "  +0.17%  "
# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.346"
$ws.Range("E51").Value = "  +0.75%  "
